$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append 5 new device rows (157-161) for regcntr_id 10002, mirroring the
# existing rows' shape: regcntr_id, device_id, lang_code, is_active, cr_by,
# cr_dtimes, eff_dtimes.
$regcntrId = 10002
$deviceIds = @(3000176, 3000177, 3000178, 3000179, 3000180)
$langCode = "eng"
$crBy = "superadmin"
$dtimes = "now()"

$row = 157
foreach ($deviceId in $deviceIds) {
    $ws.Cells.Item($row, 1).Value = $regcntrId
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = $langCode
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = $crBy
    $ws.Cells.Item($row, 6).Value = $dtimes
    $ws.Cells.Item($row, 7).Value = $dtimes
    $row++
}

# Leave the selection where the author ended up after typing the last row.
$ws.Range("B157").Select()
